$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for the "R10" rule row (8..11 AM) from
# "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the last active cell/selection as left by the edit
$ws.Range("E8").Select()
